# Update symbol list values to reflect the latest scraped crypto data.
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h); both stored as text (inlineStr).
# For numeric-looking values (columns D and E) we must force a text number
# format before assignment so Excel keeps them as text rather than silently
# converting them to numbers/percentages, then restore the "Normal" style so
# no extra cell-level style survives the edit.

function Set-TextCell {
    param(
        $ws,
        [string]$CellRef,
        [string]$NewValue
    )
    $cell = $ws.Range($CellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $NewValue
    $cell.Style = "Normal"
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @("D2", "258.13"),
    @("E2", "0.92%"),
    @("D3", "27.15"),
    @("E3", "-3.57%"),
    @("D4", "4.910"),
    @("E4", "-8.23%"),
    @("D5", "0.05958"),
    @("E5", "2.42%"),
    @("D6", "6.688"),
    @("E6", "-0.28%"),
    @("D7", "0.8708"),
    @("E7", "-0.03%"),
    @("D8", "0.9627"),
    @("E8", "7.52%"),
    @("B9", "WazirX"),
    @("C9", "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"),
    @("D9", "0.1415"),
    @("E9", "-0.06%"),
    @("B10", "LiechtensteinCryptoassetsExchange"),
    @("C10", "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"),
    @("D10", "0.03579"),
    @("E10", "3.69%"),
    @("D11", "0.07204"),
    @("E11", "-0.47%"),
    @("D12", "0.03137"),
    @("E12", "-1.33%"),
    @("D13", "0.09238"),
    @("E13", "-0.10%"),
    @("D14", "0.001550"),
    @("E14", "0.69%"),
    @("B15", "One"),
    @("C15", "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"),
    @("D15", "0.0006057"),
    @("E15", "0.01%"),
    @("B16", "TigerCash"),
    @("C16", "https://coinranking.com/coin/6hIn06L2+tigercash-tch"),
    @("D16", "0.005998"),
    @("E16", "-0.54%"),
    @("B17", "LEO"),
    @("C17", "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"),
    @("D17", "3.486"),
    @("E17", "-0.43%"),
    @("B18", "GateToken"),
    @("C18", "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"),
    @("D18", "3.221"),
    @("E18", "-0.82%"),
    @("B19", "BTSEToken"),
    @("C19", "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"),
    @("D19", "2.219"),
    @("E19", "-2.32%"),
    @("B20", "BitpandaEcosystemToken"),
    @("C20", "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"),
    @("D20", "0.3145"),
    @("E20", "-0.70%"),
    @("E21", "-0.70%"),
    @("D22", "3.528"),
    @("E22", "0.17%"),
    @("D23", "0.04248"),
    @("E23", "2.05%"),
    @("D24", "0.1379"),
    @("E24", "0.08%"),
    @("E25", "0.15%"),
    @("D26", "0.004522"),
    @("E26", "-7.24%"),
    @("E27", "0.02%"),
    @("D28", "0.0001492"),
    @("E28", "2.60%"),
    @("D40", "0.03838"),
    @("E40", "-0.34%"),
    @("D41", "0.005889"),
    @("E41", "2.30%"),
    @("D42", "0.1104"),
    @("E42", "0.36%"),
    @("D43", "0.002299"),
    @("E43", "4.58%"),
    @("D44", "0.01048"),
    @("E44", "5.62%"),
    @("D45", "0.00005491"),
    @("E45", "3.92%"),
    @("E46", "0.01%"),
    @("D47", "0.1090"),
    @("E47", "28.35%"),
    @("D48", "0.002158"),
    @("E48", "0.83%"),
    @("D49", "0.00002099"),
    @("E49", "0.01%"),
    @("D50", "0.0001999"),
    @("E50", "0.01%")
)

foreach ($change in $changes) {
    Set-TextCell $ws $change[0] $change[1]
}
